# Adjust temperature data and nutrition
# Applies the data edits made to the "block_temperature" sheet of
# state_transition.xlsx:
#   - heat_capacity values bumped for magma_block / cooled_magma_block,
#     powder_snow / snow_block, and thin_ice
#   - the "solid" state block for the minecraft:ice row corrected from
#     minecraft:packed_ice to minecraft:ice

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("block_temperature")

# heat_capacity (column B) adjustments
$ws.Range("B2").Value = 200   # minecraft:magma_block            50 -> 200
$ws.Range("B3").Value = 200   # frostedheart:cooled_magma_block  50 -> 200
$ws.Range("B4").Value = 50    # minecraft:powder_snow             5 -> 50
$ws.Range("B5").Value = 50    # minecraft:snow_block               5 -> 50
$ws.Range("B6").Value = 50    # frostedheart:thin_ice             10 -> 50

# solid-state block reference (column E) fix for the minecraft:ice row
$ws.Range("E7").Value = "minecraft:ice"

# Match the author's final cell selection on this sheet
$ws.Range("C7").Select()
